$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '29.444.76'
$ws.Range('E2').Value = '  +0.33%  '
Set-TextValue $ws.Range('D3') '1.838.38'
$ws.Range('E3').Value = '  -0.64%  '
Set-TextValue $ws.Range('D4') '0.9989'
$ws.Range('E4').Value = '  -1.39%  '
Set-TextValue $ws.Range('D5') '243.29'
$ws.Range('E5').Value = '  -0.47%  '
Set-TextValue $ws.Range('D6') '0.6263'
$ws.Range('E6').Value = '  +1.12%  '
Set-TextValue $ws.Range('D7') '0.9997'
$ws.Range('E7').Value = '  -1.18%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range('D8') '0.2961'
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range('D9') '0.07414'
$ws.Range('E9').Value = '  -0.72%  '
Set-TextValue $ws.Range('D10') '23.32'
$ws.Range('E10').Value = '  +1.09%  '
Set-TextValue $ws.Range('D11') '0.07640'
$ws.Range('E11').Value = '  -1.31%  '
Set-TextValue $ws.Range('D12') '1.833.59'
$ws.Range('E12').Value = '  -0.60%  '
Set-TextValue $ws.Range('D13') '5.010'
$ws.Range('E13').Value = '  -0.18%  '
Set-TextValue $ws.Range('D14') '0.6753'
$ws.Range('E14').Value = '  -0.11%  '
Set-TextValue $ws.Range('D15') '83.27'
$ws.Range('E15').Value = '  +0.08%  '
Set-TextValue $ws.Range('D16') '0.000009391'
$ws.Range('E16').Value = '  +3.18%  '
Set-TextValue $ws.Range('D17') '5.886'
$ws.Range('E17').Value = '  -0.46%  '
Set-TextValue $ws.Range('D18') '29.405.94'
$ws.Range('E18').Value = '  +0.36%  '
Set-TextValue $ws.Range('D19') '2.085.56'
$ws.Range('E19').Value = '  -0.21%  '
Set-TextValue $ws.Range('D20') '237.77'
$ws.Range('E20').Value = '  -0.46%  '
Set-TextValue $ws.Range('D21') '12.52'
$ws.Range('E21').Value = '  -1.36%  '
Set-TextValue $ws.Range('D22') '1.0000'
$ws.Range('E22').Value = '  -1.41%  '
Set-TextValue $ws.Range('D23') '7.356'
$ws.Range('E23').Value = '  +2.06%  '
Set-TextValue $ws.Range('D24') '1.002'
$ws.Range('E24').Value = '  -1.30%  '
Set-TextValue $ws.Range('D25') '158.77'
$ws.Range('E25').Value = '  -0.99%  '
Set-TextValue $ws.Range('D26') '0.1418'
$ws.Range('E26').Value = '  -1.34%  '
Set-TextValue $ws.Range('D27') '8.477'
$ws.Range('E27').Value = '  -0.81%  '
Set-TextValue $ws.Range('D28') '17.75'
$ws.Range('E28').Value = '  -1.09%  '
Set-TextValue $ws.Range('D29') '0.06096'
$ws.Range('E29').Value = '  +8.42%  '
Set-TextValue $ws.Range('D30') '1.495'
$ws.Range('E30').Value = '  -0.60%  '
Set-TextValue $ws.Range('D31') '1.230'
$ws.Range('E31').Value = '  +0.50%  '
Set-TextValue $ws.Range('D32') '4.091'
$ws.Range('E32').Value = '  -0.73%  '
Set-TextValue $ws.Range('D33') '4.110'
$ws.Range('E33').Value = '  -1.31%  '
Set-TextValue $ws.Range('D34') '1.863'
$ws.Range('E34').Value = '  +0.27%  '
Set-TextValue $ws.Range('D35') '1.142'
$ws.Range('E35').Value = '  -0.40%  '
Set-TextValue $ws.Range('D36') '0.7255'
$ws.Range('E36').Value = '  -3.03%  '
Set-TextValue $ws.Range('D37') '2.610'
$ws.Range('E37').Value = '  -2.12%  '
Set-TextValue $ws.Range('D38') '2.879'
$ws.Range('E38').Value = '  +1.44%  '
Set-TextValue $ws.Range('D39') '1.219.12'
$ws.Range('E39').Value = '  +0.07%  '
Set-TextValue $ws.Range('D40') '0.01761'
$ws.Range('E40').Value = '  -1.46%  '
Set-TextValue $ws.Range('D41') '6.306'
$ws.Range('E41').Value = '  -2.93%  '
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('E43').Value = '  -1.19%  '
Set-TextValue $ws.Range('D44') '1.999.33'
$ws.Range('E44').Value = '  +0.41%  '
Set-TextValue $ws.Range('D45') '101.91'
$ws.Range('E45').Value = '  +0.22%  '
Set-TextValue $ws.Range('D46') '65.54'
$ws.Range('E46').Value = '  +0.09%  '
Set-TextValue $ws.Range('D47') '0.5072'
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D48') '0.00000000121'
$ws.Range('E48').Value = '  -2.05%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D49') '9.213'
$ws.Range('E49').Value = '  +0.23%  '
Set-TextValue $ws.Range('D50') '0.4054'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('E51').Value = '  +2.19%  '
